$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1322.6666
$ws.Range("J17").Value = 1322.6666
$ws.Range("L17").Value = 3967.9998
$ws.Range("N17").Value = -4303.9998

$ws.Range("H40").Value = 1927.28
$ws.Range("I40").Value = 1321.5385
$ws.Range("J40").Value = 2583.5
$ws.Range("K40").Value = 1321.5385
$ws.Range("L40").Value = 2583.5
$ws.Range("M40").Value = -1146.5385
$ws.Range("N40").Value = -2933.5

$ws.Range("H43").Value = 2090.0833
$ws.Range("I43").Value = 2711.5715
$ws.Range("J43").Value = 1220
$ws.Range("K43").Value = 2711.5715
$ws.Range("L43").Value = 1220
$ws.Range("M43").Value = -2642.5715
$ws.Range("N43").Value = -1358

$ws.Range("H70").Value = 1223.1333
$ws.Range("I70").Value = 1352.9412
$ws.Range("J70").Value = 1053.3846
$ws.Range("K70").Value = 4058.8236
$ws.Range("L70").Value = 3160.1538
$ws.Range("M70").Value = -3788.8236
$ws.Range("N70").Value = -3700.1538

$ws.Range("H73").Value = 1223.1333
$ws.Range("I73").Value = 1352.9412
$ws.Range("J73").Value = 1053.3846
$ws.Range("K73").Value = 4058.8236
$ws.Range("L73").Value = 3160.1538
$ws.Range("M73").Value = -3122.8236
$ws.Range("N73").Value = -5032.1538

$ws.Range("H87").Value = 31058.666
$ws.Range("J87").Value = 31058.666
$ws.Range("L87").Value = 31058.666
$ws.Range("N87").Value = -33554.666

$ws.Range("H90").Value = 31058.666
$ws.Range("J90").Value = 31058.666
$ws.Range("L90").Value = 93175.99800000001
$ws.Range("N90").Value = -105655.998

$ws.Range("H129").Value = 2735.3845
$ws.Range("I129").Value = 8875.416999999999
$ws.Range("J129").Value = 893.375
$ws.Range("K129").Value = 26626.251
$ws.Range("L129").Value = 2680.125
$ws.Range("M129").Value = -21626.251
$ws.Range("N129").Value = -12680.125

$ws.Range("H137").Value = 1555.25
$ws.Range("I137").Value = 1477.8334
$ws.Range("K137").Value = 4433.5002
$ws.Range("M137").Value = -1883.5002

$ws.Range("H140").Value = 64211.285
$ws.Range("J140").Value = 64211.285
$ws.Range("L140").Value = 64211.285
$ws.Range("N140").Value = -74571.285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29514.416
$ws.Range("I32").Value = 7240.4
$ws.Range("J32").Value = 140884.5
$ws.Range("K32").Value = 7240.4
$ws.Range("L32").Value = 140884.5
$ws.Range("M32").Value = -6953.4
$ws.Range("N32").Value = -141458.5

$ws.Range("H45").Value = 71226.664
$ws.Range("I45").Value = 92614.63
$ws.Range("K45").Value = 92614.63
$ws.Range("M45").Value = -92237.63

$ws.Range("H110").Value = 43570280
$ws.Range("I110").Value = 47719750
$ws.Range("J110").Value = 805
$ws.Range("K110").Value = 47719750
$ws.Range("L110").Value = 805
$ws.Range("M110").Value = -47717705
$ws.Range("N110").Value = -4895

$ws.Range("H124").Value = 27476.334
$ws.Range("J124").Value = 27476.334
$ws.Range("L124").Value = 27476.334
$ws.Range("N124").Value = -37296.334

$ws.Range("H125").Value = 39403.168
$ws.Range("J125").Value = 39403.168
$ws.Range("L125").Value = 39403.168
$ws.Range("N125").Value = -49243.168

$ws.Range("H128").Value = 50695
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 50695
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 50695
$ws.Range("M128").ClearContents()
$ws.Range("N128").Value = -60655

$ws.Range("H130").Value = 34697.4
$ws.Range("J130").Value = 34697.4
$ws.Range("L130").Value = 34697.4
$ws.Range("N130").Value = -44737.4

$ws.Range("H132").Value = 13215.102
$ws.Range("I132").Value = 17002.771
$ws.Range("J132").Value = 3745.9285
$ws.Range("K132").Value = 51008.313
$ws.Range("L132").Value = 11237.7855
$ws.Range("M132").Value = -48478.313
$ws.Range("N132").Value = -16297.7855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 140185.75
$ws.Range("I86").Value = 185901
$ws.Range("J86").Value = 3040
$ws.Range("K86").Value = 185901
$ws.Range("L86").Value = 3040
$ws.Range("M86").Value = -184778
$ws.Range("N86").Value = -5286

$ws.Range("H89").Value = 140185.75
$ws.Range("I89").Value = 185901
$ws.Range("J89").Value = 3040
$ws.Range("K89").Value = 929505
$ws.Range("L89").Value = 15200
$ws.Range("M89").Value = -923889
$ws.Range("N89").Value = -26432

$ws.Range("H105").Value = 112963.336
$ws.Range("I105").Value = 78854.53999999999
$ws.Range("J105").Value = 201646.2
$ws.Range("K105").Value = 78854.53999999999
$ws.Range("L105").Value = 201646.2
$ws.Range("M105").Value = -77107.53999999999
$ws.Range("N105").Value = -205140.2

$ws.Range("H126").Value = 14486.667
$ws.Range("J126").Value = 14486.667
$ws.Range("L126").Value = 14486.667
$ws.Range("N126").Value = -24366.667

$ws.Range("H130").Value = 30864.389
$ws.Range("J130").Value = 30864.389
$ws.Range("L130").Value = 30864.389
$ws.Range("N130").Value = -40904.389

$ws.Range("H134").Value = 2069.1667
$ws.Range("I134").Value = 1389.6666
$ws.Range("K134").Value = 4168.9998
$ws.Range("M134").Value = -1633.9998

$ws.Range("H140").Value = 43628.57
$ws.Range("J140").Value = 43628.57
$ws.Range("L140").Value = 43628.57
$ws.Range("N140").Value = -53988.57

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 40233.4
$ws.Range("J64").Value = 40233.4
$ws.Range("L64").Value = 40233.4
$ws.Range("N64").Value = -40729.4

$ws.Range("H67").Value = 40233.4
$ws.Range("J67").Value = 40233.4
$ws.Range("L67").Value = 40233.4
$ws.Range("N67").Value = -41949.4

$ws.Range("H100").Value = 84800
$ws.Range("J100").Value = 84800
$ws.Range("L100").Value = 84800
$ws.Range("N100").Value = -86964

$ws.Range("H107").Value = 774.8095
$ws.Range("I107").Value = 736.3333
$ws.Range("J107").Value = 871
$ws.Range("K107").Value = 736.3333
$ws.Range("L107").Value = 871
$ws.Range("M107").Value = 1183.6667
$ws.Range("N107").Value = -4711

$ws.Range("H124").Value = 28830
$ws.Range("J124").Value = 28830
$ws.Range("L124").Value = 28830
$ws.Range("N124").Value = -33740

$ws.Range("H134").Value = 1883.2
$ws.Range("I134").Value = 1684.0555
$ws.Range("J134").Value = 2181.9167
$ws.Range("K134").Value = 5052.166499999999
$ws.Range("L134").Value = 6545.750100000001
$ws.Range("M134").Value = -2517.166499999999
$ws.Range("N134").Value = -11615.7501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1189.439
$ws.Range("I5").Value = 1121.6428
$ws.Range("K5").Value = 3364.9284
$ws.Range("M5").Value = -3252.9284

$ws.Range("H122").Value = 477.75
$ws.Range("I122").Value = 402.63635
$ws.Range("K122").Value = 3623.72715
$ws.Range("M122").Value = -1173.72715

$ws.Range("H135").Value = 1189.439
$ws.Range("I135").Value = 1121.6428
$ws.Range("K135").Value = 10094.7852
$ws.Range("M135").Value = -7559.7852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H127").Value = 32574
$ws.Range("J127").Value = 35000
$ws.Range("L127").Value = 35000
$ws.Range("N127").Value = -44920

$ws.Range("H132").Value = 2646.0488
$ws.Range("I132").Value = 2059.0938
$ws.Range("K132").Value = 6177.2814
$ws.Range("M132").Value = -3647.2814

$ws.Range("H136").Value = 16881.45
$ws.Range("J136").Value = 16881.45
$ws.Range("L136").Value = 50644.35000000001
$ws.Range("N136").Value = -55744.35000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2565.158
$ws.Range("J122").Value = 2888
$ws.Range("L122").Value = 8664
$ws.Range("N122").Value = -13564

$ws.Range("H127").Value = 38917.168
$ws.Range("J127").Value = 38917.168
$ws.Range("L127").Value = 38917.168
$ws.Range("N127").Value = -48837.168

$ws.Range("H128").Value = 39437.5
$ws.Range("J128").Value = 39437.5
$ws.Range("L128").Value = 39437.5
$ws.Range("N128").Value = -49397.5

$ws.Range("H130").Value = 33569.918
$ws.Range("J130").Value = 33569.918
$ws.Range("L130").Value = 33569.918
$ws.Range("N130").Value = -43609.918

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 27469.75
$ws.Range("J124").Value = 27469.75
$ws.Range("L124").Value = 27469.75
$ws.Range("N124").Value = -37289.75

$ws.Range("H126").Value = 1645.8462
$ws.Range("I126").Value = 1714
$ws.Range("J126").Value = 1492.5
$ws.Range("K126").Value = 5142
$ws.Range("L126").Value = 4477.5
$ws.Range("M126").Value = -2672
$ws.Range("N126").Value = -9417.5

$ws.Range("H131").Value = 48259
$ws.Range("J131").Value = 48259
$ws.Range("L131").Value = 48259
$ws.Range("N131").Value = -58339

$ws.Range("H141").Value = 39265
$ws.Range("J141").Value = 39265
$ws.Range("L141").Value = 39265
$ws.Range("N141").Value = -49625
